$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# row -> [B, C, D, E, F] (only columns present get updated; $null means "leave as is")
$updates = @{
    2  = @($null,               0.02839422225952148, 176092, 2484,   2516)
    3  = @($null,               0.04664087295532227, 214404, 5008,   4992)
    4  = @($null,               0.4623169898986816,  216302, 50095,  49905)
    5  = @($null,               2.476893186569214,   201866, 249836, 250164)
    6  = @($null,               6.946522235870361,   143957, 499989, 500011)
    7  = @($null,               23.27093625068665,   171888, 1998423, 2001577)
    8  = @(100000,              0.03142237663269043, 159122, 2491,   2509)
    9  = @(100000,              0.04271769523620605, 234095, 5035,   4965)
    10 = @(100000,              0.4783587455749512,  209048, 49916,  50084)
    11 = @(100000,              2.330848932266235,   214514, 249956, 250044)
    12 = @(100000,              4.790394067764282,   208751, 499771, 500229)
    13 = @(100000,              18.23332953453064,   219378, 2000921, 1999079)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $b = $vals[0]
    if ($null -ne $b) {
        $ws.Cells.Item($row, 2).Value = $b
    }

    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
}

$wb.Save()
